# Apply the "Pin data" workbook changes:
#  - R2 pin value changes from 12 to 13
#  - New "Bluetooth" section added (rows 28-30) with header fill style
#  - New shared strings: Bluetooth, 0RX, 1TX, HR, HT

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R2 value update (B20: 12 -> 13)
$ws.Range("B20").Value = 13

# New Bluetooth section
$ws.Range("A28").Value = "Bluetooth"

# Write remaining new unique strings in the same order they were first
# introduced in the saved workbook (0RX, 1TX, HR, HT) so that the shared
# string table indices line up with the target file.
$ws.Range("B29").Value = "0RX"
$ws.Range("B30").Value = "1TX"
$ws.Range("A30").Value = "HR"
$ws.Range("A29").Value = "HT"

# Give the new section header (A28) the same subtle grey shading used for
# other header rows in the sheet (White, Background 1, Darker 5%).
$ws.Range("A28").Interior.Color = 15921906

Write-Host "Applied Pin data updates."
